$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ratios" sheet lists the repeated analyses (odd rows, no shading) right
# below each original analysis (even rows, light-green shading). To make the
# extra/repeated-analysis context easier to read, every value row (2-16) is
# switched to the Arial font and the numbers are centered, for both the
# shaded and the un-shaded rows.
$rng = $ws.Range("A2:U16")
$rng.Font.Name = "Arial"
$rng.HorizontalAlignment = -4108   # xlCenter

$wb.Save()
